# Apply reviewer-requested text corrections to the response burden scheme sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: "sub-question" -> "sub-questions"
$ws.Range("A13").Value = "Answer to sub-questions of up to 2 lines"

# Row 14: fix typo "Rsponse" -> "Response"
$ws.Range("A14").Value = "a) Response to half-open question with less than 8 possibilities"

# Move the active selection to A28, matching the reviewer's final cursor position.
$ws.Range("A28").Select()
